$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K with header "Identificación médico", matching the style
# of the adjacent header cell I1 (same fill/border/bold formatting).
$ws.Range("I1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Identificación médico"

# Size the new column to fit the longer header text.
$ws.Columns("K").ColumnWidth = 18

# The header row now wraps across fewer, narrower columns -> shorter row.
$ws.Rows(1).RowHeight = 52.8

# Move the active selection to K3, as in the edited workbook.
$ws.Range("K3").Select()
